$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gene")

# Swap the start/end coordinate values in row 2 (A2 <-> B2)
$ws.Range("A2").Value = 214725146
$ws.Range("B2").Value = 214810183

# Update the active selection on the sheet to match the saved cursor position
$ws.Activate()
$ws.Range("B10").Select()
